$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pays_2")
$ws.Activate()

# --- Fix the "Afrique-Équatoriale française" wording in the "Ancien nom" column
#     (Tchad, Centrafrique, Gabon) and drop the thin-border formatting those
#     three cells used to carry, matching the now-unbordered replacement text.
$fixed = "Afrique [" + [char]0x00C9 + "|E]quatoriale française"
foreach ($addr in @("B17", "B18", "B19")) {
    $cell = $ws.Range($addr)
    $cell.Borders.LineStyle = -4142
    $cell.Value = $fixed
}

# --- India: split the single "Inde" cell into two shared-string variants
$ws.Range("A47").Value = "__Inde"
$ws.Range("B47").Value = "Inde "

# --- Add the two Congo entries (Congo-Kinshasa / Congo-Leopoldville and
#     Congo-Brazzaville) as new rows at the bottom of the table.
$ws.Range("A56").Value = "Congo belge"
$ws.Range("B56").Value = "Congo-Kinshasa|Congo-Léopoldville"
$ws.Range("C56").Value = "Belgique"
$ws.Range("D56").Value = "Afrique"
$ws.Range("E56").Borders.LineStyle = 1
$ws.Range("E56").NumberFormat = "DD/MM/YY"
$ws.Range("F56").Value = 1960

$ws.Range("A57").Value = "Congo-Brazzaville"
$ws.Range("B57").Value = $fixed
$ws.Range("C57").Value = "France"
$ws.Range("D57").Value = "Afrique"
$ws.Range("E57").Borders.LineStyle = 1
$ws.Range("E57").NumberFormat = "DD/MM/YY"
$ws.Range("F57").Value = 1960

# --- Move the selection/scroll position to where editing finished
$ws.Range("F58").Select()
